$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# Resize/reposition the subtitle placeholder textbox (EMU -> points, 1pt = 12700 EMU)
$shp.Left = 611560 / 12700
$shp.Top = 3501008 / 12700
$shp.Width = 8280920 / 12700
$shp.Height = 1752600 / 12700

# Split the first run's text into "Bài 6" + ". Mô hình ngôn ngữ"
$tr = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$run1 = $para1.Runs(1, 1)
$run1.Text = "Bài 6. Mô hình ngôn ngữ"

$newRun = $run1.Characters(1, 5)
$newRun.Text = "Bài 6"

$restRun = $run1.Characters(6, $run1.Text.Length - 5)
$restRun.Text = ". Mô hình ngôn ngữ"
